$d = $word.ActiveDocument

# The last paragraph in the document holds the "_GoBack" bookmark.
# We need to:
#   1. Insert two new empty paragraphs right before it (matching the
#      existing "en-GB" language formatting used throughout the doc).
#   2. Add the text "Updating for the second time" as a run inside that
#      same (now-last) paragraph, right before the bookmark markers.
#
# Using Range.InsertXML lets us specify the exact resulting markup for
# that paragraph (and the two new ones before it) in one shot, without
# leaving behind the stray empty runs that InsertParagraphBefore/After
# tend to introduce when duplicating the inherited run formatting.

$last = $d.Paragraphs($d.Paragraphs.Count)
$r = $last.Range

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$emptyPara = "<w:p $wNs><w:pPr><w:rPr><w:lang w:val=""en-GB""/></w:rPr></w:pPr></w:p>"

$finalPara = "<w:p $wNs w:rsidR=""00033192"" w:rsidRPr=""00AA7CBA"" w:rsidRDefault=""00033192"">" +
             "<w:pPr><w:rPr><w:lang w:val=""en-GB""/></w:rPr></w:pPr>" +
             "<w:r><w:rPr><w:lang w:val=""en-GB""/></w:rPr><w:t>Updating for the second time</w:t></w:r>" +
             "<w:bookmarkStart w:id=""0"" w:name=""_GoBack""/><w:bookmarkEnd w:id=""0""/>" +
             "</w:p>"

$xml = $emptyPara + $emptyPara + $finalPara

$r.InsertXML($xml)
